$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Bruselas (repollito)" at
# Vega Central Mapocho de Santiago. Insert a new row at position 48,
# pushing all the existing rows (48-103) down by one (to 49-104), and
# populate the newly inserted row 48 with the new record's data.
$ws.Rows("48").Insert()

$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 45128
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 100112035
$ws.Range("G48").Value = "Bruselas (repollito)"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 34
$ws.Range("K48").Value = 17000
$ws.Range("L48").Value = 17000
$ws.Range("M48").Value = 17000
$ws.Range("N48").Value = "$/malla 15 kilos"
$ws.Range("O48").Value = "Provincia de Quillota"
$ws.Range("P48").Value = 1133
$ws.Range("Q48").Value = 15
$ws.Range("R48").Value = "Hortaliza"
